# Change the "streetNumber" field's type from "int" to "string" in both
# places it appears in the document (two copies of the Persons schema
# description). Word originally stores the field description as a single
# run "streetNumber" + ": int " - editing just the word "int" (selecting
# " int" and typing " string") causes Word to split the surrounding run
# and leaves the "_GoBack" bookmark marking the last edited spot. That
# bookmark previously sat in front of "zipCode" in the second copy; after
# this edit it moves to the newly-edited first "streetNumber" occurrence.

$d = $word.ActiveDocument

# Locate every "streetNumber" occurrence (position right after the word)
# before any edits happen, searching the whole document content.
$search = $d.Content
$search.Collapse(1)   # wdCollapseStart

$occurrences = @()
while ($search.Find.Execute("streetNumber", $false, $false, $false, $false,
                             $false, $true, 0, $false, $null, 0)) {
    $occurrences += $search.End
    $search.Collapse(0)   # wdCollapseEnd - continue searching after this hit
}

# Edit from the last occurrence back to the first so earlier offsets stay
# valid while later ones are rewritten.
for ($i = $occurrences.Count - 1; $i -ge 0; $i--) {
    $afterWord = $occurrences[$i]

    # Text right after "streetNumber" is ": int " -> the selectable
    # " int" token (leading space, no trailing space) starts right after
    # the colon and is 4 characters long.
    $colonEnd = $afterWord + 1
    $typeRange = $d.Range($colonEnd, $colonEnd + 4)
    if ($typeRange.Text -ne " int") {
        throw "Unexpected text near streetNumber: [$($typeRange.Text)]"
    }

    # Replace " int" with " string", as if it had been selected and retyped.
    $typeRange.Text = " string"

    # Re-create the "_GoBack" bookmark at the boundary right after the
    # colon (this also forces/keeps the run split between ":" and
    # " string"), then immediately move it to the boundary right after
    # " string" (between " string" and the trailing space) - its real
    # final resting place. Because bookmark names are unique, adding
    # "_GoBack" again relocates it and removes it from wherever it was
    # before (including the old spot in front of "zipCode").
    $splitAfterColon = $d.Range($colonEnd, $colonEnd)
    $d.Bookmarks.Add("_GoBack", $splitAfterColon) | Out-Null

    $splitAfterString = $d.Range($colonEnd + 7, $colonEnd + 7)
    $d.Bookmarks.Add("_GoBack", $splitAfterString) | Out-Null
}
